$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the missing screw/hardware rows (18-26) ---
# Columns: A=Pedido date, B=Recibido date, C=Cantidad, D=Item, E=Precio, F=Direccion
$rows = @(
    @{ Row=18; C=4;  D="Arandelas M4";    E=0.84440000000000004; F="Tornicalvo" },
    @{ Row=19; C=90; D="Arandelas M3";    E=0.84440000000000004; F="Tornicalvo" },
    @{ Row=20; C=4;  D="Tuercas M4";      E=0.84440000000000004; F="Tornicalvo" },
    @{ Row=21; C=90; D="Tuercas M3";      E=0.84440000000000004; F="Tornicalvo" },
    @{ Row=22; C=4;  D="Tornillos M4x20"; E=0.84440000000000004; F="Tornicalvo" },
    @{ Row=23; C=30; D="Tornillos M3x10"; E=0.84440000000000004; F="Tornicalvo" },
    @{ Row=24; C=30; D="Tornillos M3x25"; E=0.84440000000000004; F="Tornicalvo" },
    @{ Row=25; C=30; D="Tornillos M3x40"; E=0.84440000000000004; F="Tornicalvo" },
    @{ Row=26; C=4;  D="Prisionero M3";   E=0.84440000000000004; F="v" }
)

foreach ($r in $rows) {
    $n = $r.Row

    # Copy the existing date formatting (style index 1) onto the new A/B cells
    $ws.Range("A6:B6").Copy() | Out-Null
    $ws.Range("A$n`:B$n").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Range("A$n").Value = 40947
    $ws.Range("B$n").Value = 40947

    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E

    # Highlight the quantity / item / price cells in yellow
    $ws.Range("C$n`:E$n").Interior.Color = 65535

    $ws.Range("F$n").Value = $r.F
}

$ws.Application.CutCopyMode = $false

# --- Add Santiago's comment on E26 ---
$commentCell = $ws.Range("E26")
$commentText = "Santiago López Pina:" + [char]10 + "No teno precio unitario de cada cosa, en el ticket de todo esto marca 7,60"
$comment = $commentCell.AddComment($commentText)
$comment.Author = "Santiago López Pina"

# --- Update selection to match the saved workbook state ---
$ws.Range("E10").Select()
